$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 444; everything from row 444 down shifts to 445+.
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new data record.
$ws.Cells.Item(444, 1).Value = 10
$ws.Cells.Item(444, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(444, 3).Value = "La Araucanía"
$ws.Cells.Item(444, 4).Value = 44461
$ws.Cells.Item(444, 5).Value = 9
$ws.Cells.Item(444, 6).Value = 100112006
$ws.Cells.Item(444, 7).Value = "Repollo"
$ws.Cells.Item(444, 8).Value = "Crespo record"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 1500
$ws.Cells.Item(444, 11).Value = 800
$ws.Cells.Item(444, 12).Value = 800
$ws.Cells.Item(444, 13).Value = 800
$ws.Cells.Item(444, 14).Value = "`$/unidad"
$ws.Cells.Item(444, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(444, 16).Value = 800
$ws.Cells.Item(444, 17).Value = 1
$ws.Cells.Item(444, 18).Value = "Hortaliza"
